$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("day_today", "Today"),
    @("day_monday", "Monday"),
    @("day_tuesday", "Tuesday"),
    @("day_wednesday", "Wednesday"),
    @("day_thursday", "Thursday"),
    @("day_friday", "Friday"),
    @("day_saturday", "Saturday"),
    @("day_sunday", "Sunday")
)

$startRow = 32
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $key = $rows[$i][0]
    $value = $rows[$i][1]

    $ws.Range("A$r").Value = $key
    $ws.Range("B$r").Value = $value
    $ws.Range("B$r").WrapText = $true
}

$ws.Range("A39").Select()
